$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.645.83"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "'2.657.58"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'597.51"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'156.71"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.647"
$ws.Range("E7").Value = "  +4.25%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.126"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Value = "'0.399"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "'28.52"
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("D14").Value = "'0.0000196"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "'3.128.40"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "'65.500.90"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'2.685.53"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "'4.76"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "'349.52"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'69.48"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "'1.78"
$ws.Range("E24").Value = "  +8.15%  "
$ws.Range("D25").Value = "'0.0000111"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'9.53"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "'1.61"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").Value = "'561.09"
$ws.Range("E28").Value = "  +5.81%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'8.06"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.162"
$ws.Range("E30").Value = "  -3.24%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").Value = "'1.80"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").Value = "'6.53"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'5.46"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").Value = "'0.420"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "'20.40"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").Value = "'155.22"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'160.40"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "'0.0604"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("D46").Value = "'22.58"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "'0.638"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").Value = "'0.102"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").Value = "'19.61"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  +5.57%  "
